# "ref: first results first point"
# The workbook's first sheet ("Datos punto 1") held placeholder/dummy test
# numbers (1,2,3,4 style sequences) in the "Volumen de datos" block
# (columns X:AY) for row 4. This edit replaces them with the first real
# measured results (C1/C2 only -- C3/C4 for this point weren't measured,
# so those cells are cleared outright), fills in the "Costos"/"Tiempos"
# C1/C2 cells (D4,E4,H4,I4) that were never populated, and formats a
# couple of the resulting rate cells with decimal number formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos punto 1")

# --- Costos / Tiempos (Explain plan block) : C1/C2 now populated -------
$ws.Range("D4").Value = 13
$ws.Range("E4").Value = 13
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1

# --- Volumen de datos block: replace dummy 1..8 placeholders -----------
# f (X:AA), g (AB:AE), h (AF:AI), i (AJ:AM), j (AN:AQ), k (AR:AU)
# Only C1/C2 were actually measured for this point; C3/C4 are cleared.
$ws.Range("X4").Value = 362
$ws.Range("Y4").Value = 362
$ws.Range("Z4").Clear()
$ws.Range("AA4").Clear()

$ws.Range("AB4").Value = 0
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Clear()
$ws.Range("AE4").Clear()

$ws.Range("AF4").Value = 5000
$ws.Range("AG4").Value = 5000
$ws.Range("AH4").Clear()
$ws.Range("AI4").Clear()

$ws.Range("AJ4").Value = 5000
$ws.Range("AK4").Value = 5000
$ws.Range("AL4").Clear()
$ws.Range("AM4").Clear()

$ws.Range("AN4").Value = 335
$ws.Range("AO4").Value = 335
$ws.Range("AP4").Clear()
$ws.Range("AQ4").Clear()

$ws.Range("AR4").Value = 0
$ws.Range("AS4").Value = 27
$ws.Range("AT4").Clear()
$ws.Range("AU4").Clear()

# --- Number formats on a few of the recomputed rate cells --------------
$ws.Range("H4").NumberFormat = "0.00"
$ws.Range("T4").NumberFormat = "0.00"
$ws.Range("U4").NumberFormat = "0.000"

# --- Selection left where the user was working after entering the data -
$ws.Range("R16").Select() | Out-Null
